$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 116; all rows from 116..159 shift down to 117..160
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly record
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = "Femacal de La Calera"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44875
$ws.Range("E116").Value = 5
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100107
$ws.Range("H116").Value = "Otros"
$ws.Range("I116").Value = 100107011
$ws.Range("J116").Value = "Tuna"
$ws.Range("K116").Value = "Sin especificar"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 70
$ws.Range("N116").Value = 32000
$ws.Range("O116").Value = 32000
$ws.Range("P116").Value = 32000
$ws.Range("Q116").Value = "`$/caja 16 kilos"
$ws.Range("R116").Value = "Provincia de Los Andes"
$ws.Range("S116").Value = 2000
$ws.Range("T116").Value = 16
